$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Black Amber, Primera) was reported and needs to
# be inserted as row 30, pushing the existing rows 30-45 down to 31-46.
$ws.Rows(30).Insert()

$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44588
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 220
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 9500
$ws.Range("P30").Value = 9227
$ws.Range("Q30").Value = "`$/caja 16 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 577
$ws.Range("T30").Value = 16
